# Add a new bulleted "Objective" line to the "Optimization settings" box
# (slide 1, shape "Rectangle 8") right after "Type of pipe model".
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(6)

$textRange = $shape.TextFrame.TextRange
$textRange.InsertAfter("`rObjective")
